$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report date header
$ws.Range("I1").Value = "18/03/2023"

# Row data: row => @(B, C, D, E, F, G, H, I, J)
$data = @{
    2  = @(2,   392, 446, 13, 0, 10, 40, 594.1, -24.92846322167985)
    3  = @(0,   27,  29,  1,  0, 0,  0,  60,    -51.66666666666666)
    4  = @(0,   84,  86,  1,  0, 5,  0,  69,    24.63768115942029)
    5  = @(2,   234, 270, 12, 0, 4,  23, 471,   -42.67515923566879)
    6  = @(24,  135, 197, 35, 3, 3,  0,  314,   -37.26114649681529)
    7  = @(1,   62,  66,  2,  1, 0,  0,  106,   -37.73584905660378)
    8  = @(0,   60,  78,  18, 0, 0,  0,  92,    -15.21739130434783)
    9  = @(0,   125, 126, 0,  0, 2,  1,  448,   -71.875)
    10 = @(0,   14,  17,  3,  0, 0,  0,  57,    -70.17543859649122)
    11 = @(0,   0,   0,   0,  0, 0,  0,  2,     -100)
    12 = @(0,   19,  21,  2,  0, 0,  0,  27,    -22.22222222222222)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
    $ws.Cells.Item($row, 9).Value = $vals[7]
    $ws.Cells.Item($row, 10).Value = $vals[8]
}
